$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 64, pushing existing rows 64:148 down to 65:149,
# mirroring the other rows' layout/style.
$ws.Rows.Item(64).Insert()

# Populate the newly-inserted row 64 with the new weekly price entry.
$ws.Cells.Item(64, 1).Value = 8
$ws.Cells.Item(64, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44721
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = 100112040
$ws.Cells.Item(64, 7).Value = "Cilantro"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 2400
$ws.Cells.Item(64, 11).Value = 1500
$ws.Cells.Item(64, 12).Value = 2000
$ws.Cells.Item(64, 13).Value = 1750
$ws.Cells.Item(64, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(64, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(64, 16).Value = 1167
$ws.Cells.Item(64, 17).Value = 1.5
$ws.Cells.Item(64, 18).Value = "Hortaliza"
